# Auto-generated COM-interop script replicating the invoice-data refresh for Sheet1.
# Row 2 is refreshed with the latest RSD totals invoice; rows 3-7 are newly
# appended invoice lines pulled in by the pipeline run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 previously carried a PO_Number (M2); the refreshed record has no PO number
# and instead populates the Job_Number column (N2), so clear the stale cell first.
$ws.Range("M2").ClearContents()

# --- Row 2 ---
foreach ($addr in @("B2","C2","H2","I2","J2","L2","P2")) { $ws.Range($addr).NumberFormat = "@" }  # keep as text (avoid date/number auto-conversion)
$ws.Range("A2").Value = "AA1"
$ws.Range("B2").Value = "08/21/25"
$ws.Range("C2").Value = "2025-08-22"
$ws.Range("D2").Value = "RSDTOT"
$ws.Range("E2").Value = "RSD"
$ws.Range("F2").Value = "I"
$ws.Range("G2").Value = "11395918-00"
$ws.Range("H2").Value = "08/21/25"
$ws.Range("I2").Value = "186.74"
$ws.Range("J2").Value = "18.13"
$ws.Range("L2").Value = "168.61"
$ws.Range("N2").Value = 25.19
$ws.Range("P2").Value = "5030"
$ws.Range("Q2").Value = 320
$ws.Range("R2").Value = "M"
$ws.Range("T2").Value = "1_1755887243708.pdf"

# --- Row 3 ---
foreach ($addr in @("B3","C3","H3","I3","J3","L3","P3")) { $ws.Range($addr).NumberFormat = "@" }  # keep as text (avoid date/number auto-conversion)
$ws.Range("A3").Value = "AA1"
$ws.Range("B3").Value = "08/21/25"
$ws.Range("C3").Value = "2025-08-22"
$ws.Range("D3").Value = "RSDTOT"
$ws.Range("E3").Value = "RSD"
$ws.Range("F3").Value = "I"
$ws.Range("G3").Value = "39395877-00"
$ws.Range("H3").Value = "08/21/25"
$ws.Range("I3").Value = "324.04"
$ws.Range("J3").Value = "29.13"
$ws.Range("L3").Value = "294.91"
$ws.Range("N3").Value = 25.16
$ws.Range("P3").Value = "5030"
$ws.Range("Q3").Value = 320
$ws.Range("R3").Value = "M"
$ws.Range("T3").Value = "2_1755887243712.pdf"

# --- Row 4 ---
foreach ($addr in @("B4","C4","H4","I4","J4","L4","P4")) { $ws.Range($addr).NumberFormat = "@" }  # keep as text (avoid date/number auto-conversion)
$ws.Range("A4").Value = "AA1"
$ws.Range("B4").Value = "08/21/25"
$ws.Range("C4").Value = "2025-08-22"
$ws.Range("D4").Value = "RSDTOT"
$ws.Range("E4").Value = "RSD"
$ws.Range("F4").Value = "I"
$ws.Range("G4").Value = "39395879-00"
$ws.Range("H4").Value = "08/21/25"
$ws.Range("I4").Value = "997.32"
$ws.Range("J4").Value = "89.64"
$ws.Range("L4").Value = "907.68"
$ws.Range("N4").Value = 24.68
$ws.Range("P4").Value = "5030"
$ws.Range("Q4").Value = 320
$ws.Range("R4").Value = "M"
$ws.Range("T4").Value = "3_1755887243715.pdf"

# --- Row 5 ---
foreach ($addr in @("B5","C5","G5","H5","I5","J5","K5","L5","P5")) { $ws.Range($addr).NumberFormat = "@" }  # keep as text (avoid date/number auto-conversion)
$ws.Range("A5").Value = "AA1"
$ws.Range("B5").Value = "08/21/25"
$ws.Range("C5").Value = "2025-08-22"
$ws.Range("D5").Value = "JONSUP"
$ws.Range("E5").Value = "Johnstone Supply"
$ws.Range("F5").Value = "I"
$ws.Range("G5").Value = "101130827.1"
$ws.Range("H5").Value = "08/21/25"
$ws.Range("I5").Value = "783.30"
$ws.Range("J5").Value = "70.40"
$ws.Range("K5").Value = "0.00"
$ws.Range("L5").Value = "712.90"
$ws.Range("O5").Value = 13616
$ws.Range("P5").Value = "5260"
$ws.Range("T5").Value = "4_1755887243717.pdf"

# --- Row 6 ---
foreach ($addr in @("B6","C6","H6","I6","J6","K6","L6","P6")) { $ws.Range($addr).NumberFormat = "@" }  # keep as text (avoid date/number auto-conversion)
$ws.Range("A6").Value = "AA1"
$ws.Range("B6").Value = "08/21/25"
$ws.Range("C6").Value = "2025-08-22"
$ws.Range("D6").Value = "LORSON"
$ws.Range("E6").Value = "Lord & Sons Inc."
$ws.Range("F6").Value = "I"
$ws.Range("H6").Value = "08/21/25"
$ws.Range("I6").Value = "31.31"
$ws.Range("J6").Value = "2.81"
$ws.Range("K6").Value = "0.00"
$ws.Range("L6").Value = "28.50"
$ws.Range("P6").Value = "1200"
$ws.Range("S6").Value = "SHOP STOCK"
$ws.Range("T6").Value = "inv-01-875854.pdf_page_1_1755887243721.pdf"

# --- Row 7 ---
foreach ($addr in @("B7","C7","H7","I7","J7","K7","L7","P7")) { $ws.Range($addr).NumberFormat = "@" }  # keep as text (avoid date/number auto-conversion)
$ws.Range("A7").Value = "AA1"
$ws.Range("B7").Value = "08/21/25"
$ws.Range("C7").Value = "2025-08-22"
$ws.Range("D7").Value = "CALHYD"
$ws.Range("E7").Value = "California Hydronics Corp"
$ws.Range("F7").Value = "I"
$ws.Range("G7").Value = "SIN221250"
$ws.Range("H7").Value = "08/21/25"
$ws.Range("I7").Value = "377.97"
$ws.Range("J7").Value = "33.97"
$ws.Range("K7").Value = "0.00"
$ws.Range("L7").Value = "344.00"
$ws.Range("M7").Value = 1504
$ws.Range("P7").Value = "5030"
$ws.Range("Q7").Value = 320
$ws.Range("R7").Value = "M"
$ws.Range("T7").Value = "sin221250_page_1_1755887243722.pdf"

# --- Column width adjustments (Excel internal width = characters + ~0.8333 padding,
#     so subtract 5/6 from the target character width to land exactly on it) ---
$ws.Columns.Item(5).ColumnWidth = 26.166666666666668   # E: 36 -> 27
$ws.Columns.Item(19).ColumnWidth = 11.166666666666666  # S: 9 -> 12
$ws.Columns.Item(20).ColumnWidth = 43.166666666666664  # T: 30 -> 44
